$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("V1").Value = "Form Tag"
